$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - first copy of the data
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 61
$wsExpo.Range("F6").Value = 5269
$wsExpo.Range("F8").Value = 5345
$wsExpo.Range("F9").Value = 620
$wsExpo.Range("F10").Value = 6
$wsExpo.Range("F11").Value = 1365

# Sheet "全部类型" (all types) - combined data, rows offset by +1 vs 展览
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 61
$wsAll.Range("F7").Value = 5269
$wsAll.Range("F9").Value = 5345
$wsAll.Range("F10").Value = 620
$wsAll.Range("F11").Value = 6
$wsAll.Range("F12").Value = 1365
